$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns I..N
$ws.Range("I1").Value = "平均时延(ms)"
$ws.Range("J1").Value = "时延抖动"
$ws.Range("K1").Value = "丢包率"
$ws.Range("L1").Value = "RTT"
$ws.Range("M1").Value = "IND"
$ws.Range("N1").Value = "CAT"

# Row 2
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 32
$ws.Range("M2").Value = 35
$ws.Range("N2").Value = 30

# Row 3
$ws.Range("I3").Value = 33
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 63
$ws.Range("M3").Value = 44
$ws.Range("N3").Value = 30

# Row 4
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 230
$ws.Range("M4").Value = 294
$ws.Range("N4").Value = 30

# Row 5
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 10
$ws.Range("L5").Value = 573
$ws.Range("M5").Value = 739
$ws.Range("N5").Value = 30

# Row 6
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 15
$ws.Range("L6").Value = 2176
$ws.Range("M6").Value = 3102
$ws.Range("N6").Value = 60

# Row 7 (only I, J, K are populated)
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 20

# Update the selection to match the target view state
$ws.Range("I1:N7").Select()
